$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 205-206; existing rows 205-214 shift down to 207-216
$ws.Rows("205:206").Insert()

# Row 205
$ws.Cells.Item(205, 1).Value = 6
$ws.Cells.Item(205, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(205, 3).Value = 'Metropolitana'
$ws.Cells.Item(205, 4).Value = 44706
$ws.Cells.Item(205, 5).Value = 13
$ws.Cells.Item(205, 6).Value = 100112022
$ws.Cells.Item(205, 7).Value = 'Arveja Verde'
$ws.Cells.Item(205, 8).Value = 'Perfection'
$ws.Cells.Item(205, 9).Value = 'Primera'
$ws.Cells.Item(205, 10).Value = 25
$ws.Cells.Item(205, 11).Value = 35000
$ws.Cells.Item(205, 12).Value = 35000
$ws.Cells.Item(205, 13).Value = 35000
$ws.Cells.Item(205, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(205, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(205, 16).Value = 1400
$ws.Cells.Item(205, 17).Value = 25
$ws.Cells.Item(205, 18).Value = 'Hortaliza'

# Row 206
$ws.Cells.Item(206, 1).Value = 6
$ws.Cells.Item(206, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(206, 3).Value = 'Metropolitana'
$ws.Cells.Item(206, 4).Value = 44706
$ws.Cells.Item(206, 5).Value = 13
$ws.Cells.Item(206, 6).Value = 100112022
$ws.Cells.Item(206, 7).Value = 'Arveja Verde'
$ws.Cells.Item(206, 8).Value = 'Perfection'
$ws.Cells.Item(206, 9).Value = 'Segunda'
$ws.Cells.Item(206, 10).Value = 20
$ws.Cells.Item(206, 11).Value = 33000
$ws.Cells.Item(206, 12).Value = 33000
$ws.Cells.Item(206, 13).Value = 33000
$ws.Cells.Item(206, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(206, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(206, 16).Value = 1320
$ws.Cells.Item(206, 17).Value = 25
$ws.Cells.Item(206, 18).Value = 'Hortaliza'
